$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.3054405439560526
$ws.Range("B1").Value = 0.3047320436772907
$ws.Range("A2").Value = -0.2586188896154695
$ws.Range("B2").Value = 0.25584808891436417
$ws.Range("A3").Value = -0.13457170876834468
$ws.Range("B3").Value = 0.1339783023954837
$ws.Range("A4").Value = -0.12197830240938501
$ws.Range("B4").Value = 0.12146770836844567
$ws.Range("A5").Value = -0.11546770841575338
$ws.Range("B5").Value = 0.11446963133783772
$ws.Range("A6").Value = -0.06494121271738518
$ws.Range("B6").Value = 0.06487945567702402
$ws.Range("A7").Value = -0.04487945573593066
$ws.Range("B7").Value = 0.04476996843975911
$ws.Range("A8").Value = -0.024769968499028572
$ws.Range("B8").Value = 0.024724058744264532
$ws.Range("A9").Value = -0.062375301207000966
$ws.Range("B9").Value = 0.062040888810958705
$ws.Range("A10").Value = -0.05604088886095582
$ws.Range("B10").Value = 0.05598950240001699
$ws.Range("A11").Value = -0.05148950244904782
$ws.Range("B11").Value = 0.05140706211999202
$ws.Range("A12").Value = -0.045407062170368384
$ws.Range("B12").Value = 0.04515674374590084
$ws.Range("A13").Value = -0.03915674379710321
$ws.Range("B13").Value = 0.03908833421294755
$ws.Range("A14").Value = -0.02708833426881796
$ws.Range("B14").Value = 0.027054908129760413
$ws.Range("A15").Value = -0.02105490818141309
$ws.Range("B15").Value = 0.02102856028461275
$ws.Range("A16").Value = -0.01502856033644906
$ws.Range("B16").Value = 0.01500492736479786
$ws.Range("A17").Value = -0.009004927416882857
$ws.Range("B17").Value = 0.008999999945645598
$ws.Range("A18").Value = -0.03611315410956806
$ws.Range("B18").Value = 0.036097442699254856
$ws.Range("A19").Value = -0.0270974427463333
$ws.Range("B19").Value = 0.0270143311274067
$ws.Range("A20").Value = -0.018014331174931897
$ws.Range("B20").Value = 0.018004351966313337
$ws.Range("A21").Value = -0.009004352013927353
$ws.Range("B21").Value = 0.008999999952343352
$ws.Range("A22").Value = -0.11227638019451902
$ws.Range("B22").Value = 0.11180685943261892
$ws.Range("A23").Value = -0.10280685948164958
$ws.Range("B23").Value = 0.10201754603829549
$ws.Range("A24").Value = -0.04212817684269332
$ws.Range("B24").Value = 0.04199999992520276
$ws.Range("A25").Value = -0.043528418668874025
$ws.Range("B25").Value = 0.04347843294194931
$ws.Range("A26").Value = -0.037478432990525334
$ws.Range("B26").Value = 0.037421196010441804
$ws.Range("A27").Value = -0.031421196059101764
$ws.Range("B27").Value = 0.031249583666911107
$ws.Range("A28").Value = 0.003090305928755477
$ws.Range("B28").Value = -0.0031146839482163458
$ws.Range("A29").Value = 0.01511468389475823
$ws.Range("B29").Value = -0.015129025705583388
$ws.Range("A30").Value = 0.035129025646302825
$ws.Range("B30").Value = -0.03534851776862169
$ws.Range("A31").Value = 0.05034851771344506
$ws.Range("B31").Value = -0.05045926580821636
$ws.Range("A32").Value = 0.06665849572747184
$ws.Range("B32").Value = -0.06683942576659518
